$d = $word.ActiveDocument

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    # Paragraph text includes trailing paragraph mark (and maybe cell mark);
    # strip control characters to check if the visible text is empty.
    $trimmed = $text -replace "[\r\a\f]", ""
    if ($trimmed -eq "") {
        $para.Range.Delete()
    }
}
